# phishing-acciones.xlsx refactor: drop two duplicate/erroneous "Ultimo
# aviso pre - Judicial" notification rows (AAREM-264967241 / AAREM-264967242,
# which mistakenly reused earlier recipients' e-mail addresses) from the
# Hoja2 log, then re-apply the AutoFilter over the cleaned-up range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# Rows 14 and 15 hold the two bogus entries; deleting row 14 twice removes
# both and shifts the rows below (16-19) up to become the new 14-17.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(14).Delete()

# Re-apply the autofilter over the new A1:E17 data extent.
$ws.Range("A1:E17").AutoFilter()

# Record the (hidden) filter-database defined name Excel keeps alongside
# an AutoFilter, scoped to this sheet.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Hoja2!`$A`$1:`$E`$17")
$filterName.Visible = $false

# Restore the saved selection.
$ws.Range("C15").Select()
